$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'youth knee sleeves with pads'
$ws.Range("A2").Value = 'youth gloves and knee pads'
$ws.Range("A3").Value = 'mens tights 3 4'
$ws.Range("A4").Value = 'black leggings capri'
$ws.Range("A5").Value = 'knee pad tights'
$ws.Range("A6").Value = 'knee leggings'
$ws.Range("A7").Value = 'mens compression pants green'
$ws.Range("A8").Value = 'men tights leggings'
$ws.Range("A9").Value = 'knee pad workout'
$ws.Range("A10").Value = 'hex knee pads basketball'
$ws.Range("A11").Value = 'white knee pads basketball'
$ws.Range("A12").Value = 'knee pad basketball leggings'
$ws.Range("A13").Value = 'capri compression tights'
$ws.Range("A14").Value = 'compression workout capri'
$ws.Range("A15").Value = 'knee pad leggings youth'
$ws.Range("A16").Value = 'all white knee pads for basketball'
$ws.Range("A17").Value = 'basketball compression pants with knee pads youth'
$ws.Range("A18").Value = 'triple 8 knee pads youth'
$ws.Range("A19").Value = 'mens leggings compression under armour'
$ws.Range("A20").Value = 'dirt bike knee pads youth'
$ws.Range("A21").Value = 'nike leggings men basketball'
$ws.Range("A22").Value = 'gym shark mens leggings'
$ws.Range("A23").Value = 'multicam pants with knee pads'
$ws.Range("A24").Value = 'tactical pants with knee pads for men'
$ws.Range("A25").Value = 'excersize gear for men'
$ws.Range("A26").Value = 'men capri leggings'
$ws.Range("A27").Value = 'mens leggings black'
$ws.Range("A28").Value = 'under armour tights'
$ws.Range("A29").Value = 'guys tights'
$ws.Range("A30").Value = 'leggings for basketball'
$ws.Range("A31").Value = 'dri fit compression pants men'
$ws.Range("A32").Value = 'pink compression pants men'
$ws.Range("A33").Value = 'men workout compression pants'
$ws.Range("A34").Value = 'mens football gear'
$ws.Range("A35").Value = 'compression pants with pads basketball'
$ws.Range("A36").Value = 'mens white leggings compression'
$ws.Range("A37").Value = 'cold gear mens'
$ws.Range("A38").Value = 'basketball leggings with kneepads'
$ws.Range("A39").Value = 'mens athletic leggins'
$ws.Range("A40").Value = 'mens compression tights basketball'
$ws.Range("A41").Value = 'under armour cold gear compression pants men'
$ws.Range("A42").Value = 'green mens compression leggings'
$ws.Range("A43").Value = 'mens compression tights leggings'
$ws.Range("A44").Value = 'compression knee pads for basketball'
$ws.Range("A45").Value = 'legging basketball men'
$ws.Range("A46").Value = 'basketball compression pants with padded knees'
$ws.Range("A47").Value = 'knee pad pants'
$ws.Range("A48").Value = 'kneepad honeycomb'
$ws.Range("A49").Value = 'black basketball knee pads'
$ws.Range("A50").Value = 'cheap knee pads for basketball'
$ws.Range("A51").Value = 'compression pants'
$ws.Range("A52").Value = 'knee pads pants'
$ws.Range("A53").Value = 'basketball youth compression pants'
$ws.Range("A54").Value = 'knee pad hex'
$ws.Range("A55").Value = 'men tights and leggings'
$ws.Range("A56").Value = 'knee pad lacrosse'
$ws.Range("A57").Value = 'mens compression tights'
$ws.Range("A58").Value = 'patella guard'
$ws.Range("A59").Value = 'squat pad knee'
$ws.Range("A60").Value = 'best knee pads basketball'
$ws.Range("A61").Value = 'kneepad basketball'
$ws.Range("A62").Value = 'compression mens running pants'
$ws.Range("A63").Value = 'knee protection soccer'
$ws.Range("A64").Value = 'sliding leg guard'
$ws.Range("A65").Value = 'athletic leggings youth'
$ws.Range("A66").Value = 'basketball tights boys'
$ws.Range("A67").Value = 'knee pads xxl'
$ws.Range("A68").Value = 'wrestling clothes for boys'
$ws.Range("A69").Value = 'boys xxl baseball pants'
$ws.Range("A70").Value = 'knee pads girls volleyball'
$ws.Range("A71").Value = 'youth boys leggings sports'
$ws.Range("A72").Value = 'mens compression running pants'
$ws.Range("A73").Value = 'black knee pads wrestling'
$ws.Range("A74").Value = 'black knee pads youth'
$ws.Range("A75").Value = 'hex pad'
$ws.Range("A76").Value = 'men spandex pants'
$ws.Range("A77").Value = 'knee protectors for men'
$ws.Range("A78").Value = 'volleyball gear'
$ws.Range("A79").Value = 'capris for men'
$ws.Range("A80").Value = 'youth black baseball pants'
$ws.Range("A81").Value = 'knee pads professional'
$ws.Range("A82").Value = 'mens mesh pants'
$ws.Range("A83").Value = 'knee pad work pants'
$ws.Range("A84").Value = 'knee pads small'
$ws.Range("A85").Value = 'youth leggings sports'
$ws.Range("A86").Value = 'fitness knee pads'
$ws.Range("A87").Value = 'knee pad sports'
$ws.Range("A88").Value = 'basketball clothes'
$ws.Range("A89").Value = 'girls workout leggings'
$ws.Range("A90").Value = 'baseball pants mens'
$ws.Range("A91").Value = 'kids basketball leggings with knee pads'
$ws.Range("A92").Value = 'youth asics wrestling knee pads'
$ws.Range("A93").Value = 'tesla wintergear for men'
$ws.Range("A94").Value = 'mcdavid basketball knee pads white'
$ws.Range("A95").Value = 'thermo ball mens'
$ws.Range("A96").Value = 'mens nike thermal training pants'
$ws.Range("A97").Value = 'man winter leggings'
$ws.Range("A98").Value = 'underware pants men'
$ws.Range("A99").Value = 'under armor youth basketball compression pants'
$ws.Range("A100").Value = 'black capris'
